# Commit: "Add files via upload"
#
# The workbook was re-uploaded after two visible, user-driven changes in
# Excel: the worksheet was renamed and the on-screen selection moved.
# Reproduce both:
#
#   1. Rename the worksheet "basic_properties_data" -> "OGS"
#   2. Move the active selection from T9 to G81
#
# (The surrounding window-chrome / co-authoring bookkeeping in the diff -
# bookViews xWindow/yWindow and xr:revisionPtr's documentId - are Excel
# client/session identifiers that Excel regenerates on every save; they
# carry no workbook content and aren't exposed on the Excel object model,
# so there is nothing for a COM script to set there.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet
$ws.Name = "OGS"

# 2) Move/scroll the selection to the new active cell
$ws.Range("G81").Select() | Out-Null
